# Swap the "M"/"T" values between columns B and C for the Friday/Saturday
# pairs of rows that make up the work-schedule pattern, and fix the
# standalone exception at C365.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where column B currently holds "M" and column C currently holds "T"
# and the two values need to be swapped (B -> "T", C -> "M").
$rows = @(11,12,18,19,25,26,32,33,39,40,46,47,53,54,60,61,67,68,74,75,81,82,88,89,95,96,102,103,109,110,116,117,123,124,130,131,137,138,144,145,151,152,158,159,165,166,172,173,179,180,186,187,193,194,200,201,207,208,214,215,221,222,228,229,235,236,242,243,249,250,256,257,263,264,270,271,277,278,284,285,291,292,298,299,305,306,312,313,319,320,326,327,333,334,340,341,347,348,354,355,361,362)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "T"
    $ws.Cells.Item($r, 3).Value = "M"
}

# Standalone fix: row 365, column C changes from "T" to "M" (column B is
# untouched there, it stays "-").
$ws.Cells.Item(365, 3).Value = "M"
